$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# ---------------------------------------------------------------------------
# Row 1 - was just "VIA_Chromium" in A1, now a full row of selector ids
# ---------------------------------------------------------------------------
Set-TextCell "A1" "VIA_Chromium"
Set-TextCell "B1" "id=startdate"
Set-TextCell "C1" "id=insurancesum"
Set-TextCell "D1" "id=meritrating"
Set-TextCell "E1" "id=damageinsurance"
Set-TextCell "F1" 'text="Euro Protection"'
Set-TextCell "G1" 'text="Legal Defense Insurance"'
Set-TextCell "H1" "id=courtesycar"
Set-TextCell "I1" "id=preventerinsurancedata"
Set-TextCell "J1" "id=nextselectpriceoption"

# ---------------------------------------------------------------------------
# Row 2 - A2 keeps its value; B2:J2 become empty cells highlighted yellow
# ---------------------------------------------------------------------------
Set-TextCell "A2" "VIA_Pixel9Pro_API35"
$ws.Range("B2:J2").NumberFormat = "@"
$ws.Range("B2:J2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Row 3 - header/label row; a few controls were renamed from combo- to
# select-prefixed ids (keep existing row style)
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "selInsuranceSum"
$ws.Range("D3").Value = "selMeritRating"
$ws.Range("E3").Value = "selDamageInsurance"
$ws.Range("H3").Value = "selCourtesyCar"

# ---------------------------------------------------------------------------
# Row 4 - "Check defaults" row gains two more default-state values
# ---------------------------------------------------------------------------
Set-TextCell "A4" "Check defaults"
Set-TextCell "F4" "<UNCHECKED>"
Set-TextCell "G4" "<UNCHECKED>"

# ---------------------------------------------------------------------------
# Row 5 - brand new smoke-test data row
# ---------------------------------------------------------------------------
Set-TextCell "A5" "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
Set-TextCell "B5" "05/01/2025"
Set-TextCell "C5" "10.000.000,00"
Set-TextCell "D5" "Bonus 5"
Set-TextCell "E5" "Partial Coverage"
Set-TextCell "F5" "<CHECK>"
Set-TextCell "G5" "<CHECK>"
Set-TextCell "H5" "Yes"
Set-TextCell "J5" "X"
